$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.902.59'
$ws.Range('E2').Value = '  +2.32%  '

$ws.Range('D3').Value = '1.810.55'
$ws.Range('E3').Value = '  +1.42%  '

$ws.Range('D4').Value = "'0.9995"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.82%  '

$ws.Range('D5').Value = "'337.25"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.64%  '

$ws.Range('D6').Value = "'0.9966"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.89%  '

$ws.Range('D7').Value = "'0.3921"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.59%  '

$ws.Range('D8').Value = "'0.3482"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.92%  '

$ws.Range('D9').Value = "'48.12"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.30%  '

$ws.Range('E10').Value = '  +0.22%  '

$ws.Range('D11').Value = "'0.07573"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.94%  '

$ws.Range('D12').Value = "'0.9970"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.85%  '

$ws.Range('D13').Value = "'22.12"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.23%  '

$ws.Range('D14').Value = "'6.523"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.27%  '

$ws.Range('D15').Value = '1.812.00'
$ws.Range('E15').Value = '  +1.31%  '

$ws.Range('D16').Value = "'7.200"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.68%  '

$ws.Range('E17').Value = '  +1.71%  '

$ws.Range('D18').Value = "'0.06688"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.59%  '

$ws.Range('D19').Value = "'85.17"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.09%  '

$ws.Range('D20').Value = "'0.9966"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.58%  '

$ws.Range('D21').Value = "'17.85"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.35%  '

$ws.Range('D22').Value = "'6.573"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.00%  '

$ws.Range('D23').Value = '27.904.00'
$ws.Range('E23').Value = '  +2.38%  '

$ws.Range('D24').Value = "'12.87"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.06%  '

$ws.Range('D25').Value = "'2.409"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.97%  '

$ws.Range('D26').Value = "'2.554"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.43%  '

$ws.Range('D27').Value = "'1.473"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.71%  '

$ws.Range('E28').Value = '  +0.23%  '

$ws.Range('D29').Value = "'154.67"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.93%  '

$ws.Range('D30').Value = '2.016.47'
$ws.Range('E30').Value = '  +1.26%  '

$ws.Range('D31').Value = "'135.46"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.96%  '

$ws.Range('E32').Value = '  -0.45%  '

$ws.Range('D33').Value = "'6.120"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.64%  '

$ws.Range('E34').Value = '  +2.41%  '

$ws.Range('D35').Value = "'13.29"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.15%  '

$ws.Range('D36').Value = "'5.540"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.60%  '

$ws.Range('D37').Value = "'0.02430"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.05%  '

$ws.Range('B38').Value = 'TheSandbox'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D38').Value = "'0.6910"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.88%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = "'0.06526"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.88%  '

$ws.Range('D40').Value = "'1.611"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.25%  '

$ws.Range('D41').Value = "'0.2225"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.79%  '

$ws.Range('D42').Value = "'1.271"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.58%  '

$ws.Range('D43').Value = "'8.565"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.38%  '

$ws.Range('D44').Value = "'14.65"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.98%  '

$ws.Range('D45').Value = "'0.6544"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.14%  '

$ws.Range('D46').Value = "'0.9968"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.66%  '

$ws.Range('D47').Value = "'3.863"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.57%  '

$ws.Range('D48').Value = "'2.162"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.46%  '

$ws.Range('D49').Value = "'132.35"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.45%  '

$ws.Range('D50').Value = "'0.07200"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.38%  '

$ws.Range('D51').Value = "'80.76"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.11%  '
